# Add 50 new accounts (rows 102-151) to the Production sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
    "williams.daft.1236@faiyamrahman.com",
    "metro.earth.1457@faiyamrahman.com",
    "fossil.fossil.1380@faiyamrahman.com",
    "metro.faiyam.597@faiyamrahman.com",
    "rahman.daft.760@faiyamrahman.com",
    "rahman.williams.697@faiyamrahman.com",
    "metro.williams.1661@faiyamrahman.com",
    "rahman.daft.1669@faiyamrahman.com",
    "rahman.williams.919@faiyamrahman.com",
    "daft.williams.1509@faiyamrahman.com",
    "bts.daft.516@faiyamrahman.com",
    "earth.water.785@faiyamrahman.com",
    "williams.grassfed.1983@faiyamrahman.com",
    "bts.grassfed.1311@faiyamrahman.com",
    "grassfed.metro.1425@faiyamrahman.com",
    "grassfed.fossil.708@faiyamrahman.com",
    "earth.faiyam.1142@faiyamrahman.com",
    "metro.bts.15@faiyamrahman.com",
    "water.grassfed.484@faiyamrahman.com",
    "rahman.bts.894@faiyamrahman.com",
    "williams.metro.831@faiyamrahman.com",
    "rahman.earth.923@faiyamrahman.com",
    "earth.rahman.1730@faiyamrahman.com",
    "fossil.bts.428@faiyamrahman.com",
    "faiyam.earth.739@faiyamrahman.com",
    "faiyam.fossil.153@faiyamrahman.com",
    "faiyam.faiyam.1624@faiyamrahman.com",
    "williams.earth.1048@faiyamrahman.com",
    "bts.williams.1249@faiyamrahman.com",
    "bts.rahman.1114@faiyamrahman.com",
    "fossil.grassfed.1308@faiyamrahman.com",
    "faiyam.williams.694@faiyamrahman.com",
    "rahman.williams.242@faiyamrahman.com",
    "bts.daft.132@faiyamrahman.com",
    "fossil.water.365@faiyamrahman.com",
    "bts.grassfed.1860@faiyamrahman.com",
    "williams.faiyam.12@faiyamrahman.com",
    "earth.bts.1735@faiyamrahman.com",
    "rahman.rahman.602@faiyamrahman.com",
    "grassfed.bts.805@faiyamrahman.com",
    "williams.rahman.1847@faiyamrahman.com",
    "bts.rahman.1030@faiyamrahman.com",
    "faiyam.bts.822@faiyamrahman.com",
    "williams.faiyam.109@faiyamrahman.com",
    "williams.bts.1129@faiyamrahman.com",
    "rahman.grassfed.1726@faiyamrahman.com",
    "grassfed.faiyam.592@faiyamrahman.com",
    "faiyam.williams.524@faiyamrahman.com",
    "fossil.earth.1755@faiyamrahman.com",
    "fossil.daft.1579@faiyamrahman.com"
)

for ($i = 0; $i -lt $emails.Count; $i++) {
    $row = 102 + $i
    $ws.Cells.Item($row, 1).Value = 100 + $i
    $ws.Cells.Item($row, 2).Value = $emails[$i]
    $ws.Cells.Item($row, 3).Value = "n/a"
    $ws.Cells.Item($row, 4).Value = "beatthestreak1"
}

# Widen column A (no longer auto best-fit) and move the selection to
# reflect the end of the newly-entered data, matching the editor state
# after the bulk entry.
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Range("B155").Select()
